$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "26.498.22"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.01%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.730.96"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("E4").Value = "  +0.07%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "243.91"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.59%  "

$ws.Range("E6").Value = "  +0.12%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4888"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +1.65%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2620"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.85%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06182"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.67%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "1.731.19"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.20%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07025"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.80%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "15.50"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.93%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.556"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.77%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.6016"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -2.22%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "77.40"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.31%  "

$ws.Range("E16").Value = "  +0.11%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "26.498.87"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.14%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.000007086"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +2.16%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "11.40"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -2.05%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "1.955.68"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.45%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.473"
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "8.599"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -4.00%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "5.184"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.80%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "138.81"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.55%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "15.27"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.52%  "

$ws.Range("E27").Value = "  +0.68%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "106.55"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.30%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.716"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -4.40%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "3.974"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.00%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.07961"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.87%  "

$ws.Range("E32").Value = "  -0.11%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.04518"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.91%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.12%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.617"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("E36").Value = "  +0.91%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.6251"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -1.81%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.9077"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -2.89%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.995"
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.413"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.55%  "

$ws.Range("E41").Value = "  +0.00%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.01489"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -0.96%  "

$ws.Range("E43").Value = "  -4.49%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "5.449"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -2.54%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.3868"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.90%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "6.679"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -3.59%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.1157"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.27%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.05366"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.66%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "30.34"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -1.79%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "7.708"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.62%  "

$ws.Range("E51").Value = "  -1.55%  "
